# Auto-generated Excel COM-interop script
# Applies numeric value updates to "Leve Profit" style sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per the target diff. Cells are plain cached numbers (no formulas in source workbook).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7554.7
$ws.Range("I40").Value = 5925.1665
$ws.Range("K40").Value = 5925.1665
$ws.Range("M40").Value = -5750.1665
$ws.Range("H43").Value = 4382
$ws.Range("J43").Value = 4500
$ws.Range("L43").Value = 4500
$ws.Range("N43").Value = -4638
$ws.Range("H70").Value = 3841.9614
$ws.Range("I70").Value = 2566.111
$ws.Range("J70").Value = 4517.4116
$ws.Range("K70").Value = 7698.333
$ws.Range("L70").Value = 13552.2348
$ws.Range("M70").Value = -7428.333
$ws.Range("N70").Value = -14092.2348
$ws.Range("H73").Value = 3841.9614
$ws.Range("I73").Value = 2566.111
$ws.Range("J73").Value = 4517.4116
$ws.Range("K73").Value = 7698.333
$ws.Range("L73").Value = 13552.2348
$ws.Range("M73").Value = -6762.333
$ws.Range("N73").Value = -15424.2348
$ws.Range("H74").Value = 10933.223
$ws.Range("I74").Value = 6679.8
$ws.Range("J74").Value = 16250
$ws.Range("K74").Value = 6679.8
$ws.Range("L74").Value = 16250
$ws.Range("M74").Value = -5743.8
$ws.Range("N74").Value = -18122
$ws.Range("H77").Value = 10933.223
$ws.Range("I77").Value = 6679.8
$ws.Range("J77").Value = 16250
$ws.Range("K77").Value = 33399
$ws.Range("L77").Value = 81250
$ws.Range("M77").Value = -28719
$ws.Range("N77").Value = -90610
$ws.Range("H107").Value = 180
$ws.Range("I107").Value = 180
$ws.Range("K107").Value = 180
$ws.Range("M107").Value = 1740
$ws.Range("H112").Value = 1833.7368
$ws.Range("J112").Value = 1622.7333
$ws.Range("L112").Value = 4868.199900000001
$ws.Range("N112").Value = -7084.199900000001
$ws.Range("H116").Value = 2779
$ws.Range("I116").Value = 1202.5
$ws.Range("K116").Value = 1202.5
$ws.Range("M116").Value = 2239.5
$ws.Range("H125").Value = 1692.4
$ws.Range("I125").Value = 1490.5
$ws.Range("K125").Value = 13414.5
$ws.Range("M125").Value = -10954.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 8933515
$ws.Range("I102").Value = 15626486
$ws.Range("J102").Value = 9553.166999999999
$ws.Range("K102").Value = 15626486
$ws.Range("L102").Value = 9553.166999999999
$ws.Range("M102").Value = -15624864
$ws.Range("N102").Value = -12797.167

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 9500
$ws.Range("I75").Value = 4000
$ws.Range("J75").Value = 15000
$ws.Range("K75").Value = 4000
$ws.Range("L75").Value = 15000
$ws.Range("M75").Value = -3064
$ws.Range("N75").Value = -16872
$ws.Range("H78").Value = 9500
$ws.Range("I78").Value = 4000
$ws.Range("J78").Value = 15000
$ws.Range("K78").Value = 12000
$ws.Range("L78").Value = 45000
$ws.Range("M78").Value = -7320
$ws.Range("N78").Value = -54360
$ws.Range("H80").Value = 905.7273
$ws.Range("I80").Value = 458.2
$ws.Range("J80").Value = 1278.6666
$ws.Range("K80").Value = 458.2
$ws.Range("L80").Value = 1278.6666
$ws.Range("M80").Value = 539.8
$ws.Range("N80").Value = -3274.6666
$ws.Range("H82").Value = 14372.363
$ws.Range("I82").Value = 5820.8
$ws.Range("K82").Value = 5820.8
$ws.Range("M82").Value = -5437.8
$ws.Range("H83").Value = 905.7273
$ws.Range("I83").Value = 458.2
$ws.Range("J83").Value = 1278.6666
$ws.Range("K83").Value = 2291
$ws.Range("L83").Value = 6393.333000000001
$ws.Range("M83").Value = 2701
$ws.Range("N83").Value = -16377.333
$ws.Range("H85").Value = 14372.363
$ws.Range("I85").Value = 5820.8
$ws.Range("K85").Value = 5820.8
$ws.Range("M85").Value = -4494.8
$ws.Range("H86").Value = 5610.5835
$ws.Range("I86").Value = 1721.5
$ws.Range("K86").Value = 1721.5
$ws.Range("M86").Value = -598.5
$ws.Range("H89").Value = 5610.5835
$ws.Range("I89").Value = 1721.5
$ws.Range("K89").Value = 8607.5
$ws.Range("M89").Value = -2991.5
$ws.Range("H99").Value = 333334000
$ws.Range("I99").Value = 333334000
$ws.Range("K99").Value = 333334000
$ws.Range("M99").Value = -333332502

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2094.625
$ws.Range("J22").Value = 2617.4
$ws.Range("L22").Value = 2617.4
$ws.Range("N22").Value = -3317.4
$ws.Range("H70").Value = 31916.666
$ws.Range("J70").Value = 31916.666
$ws.Range("L70").Value = 31916.666
$ws.Range("N70").Value = -32546.666
$ws.Range("H73").Value = 31916.666
$ws.Range("J73").Value = 31916.666
$ws.Range("L73").Value = 31916.666
$ws.Range("N73").Value = -34100.666
$ws.Range("H93").Value = 5418.625
$ws.Range("I93").Value = 1843
$ws.Range("K93").Value = 1843
$ws.Range("M93").Value = 29
$ws.Range("H99").Value = 2475.6667
$ws.Range("I99").Value = 2471.8
$ws.Range("J99").Value = 2480.5
$ws.Range("K99").Value = 2471.8
$ws.Range("L99").Value = 2480.5
$ws.Range("M99").Value = -973.8000000000002
$ws.Range("N99").Value = -5476.5
$ws.Range("H105").Value = 2554.4
$ws.Range("I105").Value = 1219
$ws.Range("K105").Value = 1219
$ws.Range("M105").Value = 528
$ws.Range("H126").Value = 2475.6667
$ws.Range("I126").Value = 2471.8
$ws.Range("J126").Value = 2480.5
$ws.Range("K126").Value = 7415.400000000001
$ws.Range("L126").Value = 7441.5
$ws.Range("M126").Value = -4945.400000000001
$ws.Range("N126").Value = -12381.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 671.25
$ws.Range("I107").Value = 325
$ws.Range("K107").Value = 975
$ws.Range("M107").Value = 945

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H131").Value = 16999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5851.8
$ws.Range("I46").Value = 4592.5713
$ws.Range("J46").Value = 7454.4546
$ws.Range("K46").Value = 4592.5713
$ws.Range("L46").Value = 7454.4546
$ws.Range("M46").Value = -4404.5713
$ws.Range("N46").Value = -7830.4546
$ws.Range("H68").Value = 4713.9165
$ws.Range("I68").Value = 1597.25
$ws.Range("J68").Value = 6272.25
$ws.Range("K68").Value = 1597.25
$ws.Range("L68").Value = 6272.25
$ws.Range("M68").Value = -848.25
$ws.Range("N68").Value = -7770.25
$ws.Range("H71").Value = 4713.9165
$ws.Range("I71").Value = 1597.25
$ws.Range("J71").Value = 6272.25
$ws.Range("K71").Value = 7986.25
$ws.Range("L71").Value = 31361.25
$ws.Range("M71").Value = -4242.25
$ws.Range("N71").Value = -38849.25
$ws.Range("H82").Value = 4414.7
$ws.Range("J82").Value = 5312.125
$ws.Range("L82").Value = 5312.125
$ws.Range("N82").Value = -6034.125
$ws.Range("H85").Value = 4414.7
$ws.Range("J85").Value = 5312.125
$ws.Range("L85").Value = 5312.125
$ws.Range("N85").Value = -7808.125
$ws.Range("H93").Value = 2622.5
$ws.Range("I93").Value = 2622.5
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2622.5
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -1374.5
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5777.4546
$ws.Range("I62").Value = 2325
$ws.Range("J62").Value = 7750.2856
$ws.Range("K62").Value = 2325
$ws.Range("L62").Value = 7750.2856
$ws.Range("M62").Value = -1701
$ws.Range("N62").Value = -8998.285599999999
$ws.Range("H65").Value = 5777.4546
$ws.Range("I65").Value = 2325
$ws.Range("J65").Value = 7750.2856
$ws.Range("K65").Value = 11625
$ws.Range("L65").Value = 38751.428
$ws.Range("M65").Value = -8505
$ws.Range("N65").Value = -44991.428
$ws.Range("H113").Value = 1092.3077
$ws.Range("I113").Value = 885.5714
$ws.Range("K113").Value = 2656.7142
$ws.Range("M113").Value = -486.7142000000003

